$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 117, shifting existing rows 117-127 down to 118-128.
$ws.Rows("117").Insert()

# Populate the newly inserted row 117 with the new weekly record.
$ws.Range("A117").Value2 = 10
$ws.Range("B117").Value2 = "Vega Modelo de Temuco"
$ws.Range("C117").Value2 = "La Araucanía"
$ws.Range("D117").Value2 = 44461
$ws.Range("E117").Value2 = 9
$ws.Range("F117").Value2 = 100112005
$ws.Range("G117").Value2 = "Puerro"
$ws.Range("H117").Value2 = "Azul de Maquehue"
$ws.Range("I117").Value2 = "Primera"
$ws.Range("J117").Value2 = 20
$ws.Range("K117").Value2 = 8000
$ws.Range("L117").Value2 = 8000
$ws.Range("M117").Value2 = 8000
$ws.Range("N117").Value2 = "`$/docena de paquetes"
$ws.Range("O117").Value2 = "Provincia de Cautín"
$ws.Range("P117").Value2 = 667
$ws.Range("Q117").Value2 = 12
$ws.Range("R117").Value2 = "Hortaliza"
